# Added Hawkeye for lbm benchmark
# Fills in the previously-blank Hawkeye (row 45) and OPTGen (row 46) rows
# for the lbm benchmark on both the Config1 and Config2 sheets, and moves
# the active selection to reflect where the author ended up working.

$wb = $excel.ActiveWorkbook

# --- Config1 sheet ---
$ws1 = $wb.Worksheets.Item("Config1")

$ws1.Range("C45").Value = 50000000
$ws1.Range("D45").Value = 82429919
$ws1.Range("E45").Value = 2641903
$ws1.Range("F45").Value = 883205
$ws1.Range("G45").Value = 1758698
$ws1.Range("H45").Formula = "=(C45/D45)"
$ws1.Range("I45").Formula = "=F45/(C45/1000)"

$ws1.Range("C46").Value = 50000000
$ws1.Range("D46").Value = 82429919
$ws1.Range("E46").Value = 49369
$ws1.Range("F46").Value = 13788
$ws1.Range("G46").Formula = "=E46-F46"
$ws1.Range("H46").Formula = "=(C46/D46)"
$ws1.Range("I46").Formula = "=G46/(C46/1000)"
$ws1.Range("J46").Formula = "=F46/E46"

[void]$ws1.Range("A47").Select()

# --- Config2 sheet ---
$ws2 = $wb.Worksheets.Item("Config2")

$ws2.Range("C45").Value = 50000000
$ws2.Range("D45").Value = 71437360
$ws2.Range("E45").Value = 2656558
$ws2.Range("F45").Value = 867959
$ws2.Range("G45").Value = 1788599
$ws2.Range("H45").Formula = "=(C45/D45)"
$ws2.Range("I45").Formula = "=G45/(C45/1000)"

$ws2.Range("C46").Value = 50000000
$ws2.Range("D46").Value = 71437360
$ws2.Range("E46").Value = 42802
$ws2.Range("F46").Value = 13968
$ws2.Range("G46").Formula = "=E46-F46"
$ws2.Range("H46").Formula = "=(C46/D46)"
$ws2.Range("I46").Formula = "=G46/(C46/1000)"
$ws2.Range("J46").Formula = "=F46/E46"

$ws2.Activate()
[void]$ws2.Range("C47").Select()

$ws1.Activate()
[void]$ws1.Range("A47").Select()
